$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.019897960467196
$ws.Cells.Item(2, 4).Value = 1.029649569611741
$ws.Cells.Item(2, 5).Value = 1.020976119449794
$ws.Cells.Item(2, 6).Value = 1.038187441699317
$ws.Cells.Item(2, 9).Value = 1.029340052793341
$ws.Cells.Item(2, 10).Value = 1.025098124277093
$ws.Cells.Item(2, 11).Value = 1.032462867119074
$ws.Cells.Item(2, 12).Value = 1.023814814373924
$ws.Cells.Item(2, 13).Value = 1.040976189024601
$ws.Cells.Item(2, 14).Value = 1.012440282311641
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.020736951672776
$ws.Cells.Item(3, 4).Value = 1.030289596716294
$ws.Cells.Item(3, 5).Value = 1.021684971002907
$ws.Cells.Item(3, 6).Value = 1.039062383187478
$ws.Cells.Item(3, 9).Value = 1.029441030009711
$ws.Cells.Item(3, 10).Value = 1.025574646747839
$ws.Cells.Item(3, 11).Value = 1.032911536577056
$ws.Cells.Item(3, 12).Value = 1.024330287406961
$ws.Cells.Item(3, 13).Value = 1.041660920558679
$ws.Cells.Item(3, 14).Value = 1.012599580062058
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.021280529922716
$ws.Cells.Item(4, 4).Value = 1.03070413859546
$ws.Cells.Item(4, 5).Value = 1.022144639558881
$ws.Cells.Item(4, 6).Value = 1.039629474038795
$ws.Cells.Item(4, 9).Value = 1.029505130081039
$ws.Cells.Item(4, 10).Value = 1.025883048816012
$ws.Cells.Item(4, 11).Value = 1.033201570173231
$ws.Cells.Item(4, 12).Value = 1.024664155521476
$ws.Cells.Item(4, 13).Value = 1.042104294687359
$ws.Cells.Item(4, 14).Value = 1.012702628241446
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.021509214954955
$ws.Cells.Item(5, 4).Value = 1.030878506290211
$ws.Cells.Item(5, 5).Value = 1.022338120418127
$ws.Cells.Item(5, 6).Value = 1.039868103335148
$ws.Cells.Item(5, 9).Value = 1.029531780474068
$ws.Cells.Item(5, 10).Value = 1.026012714182042
$ws.Cells.Item(5, 11).Value = 1.033323430563874
$ws.Cells.Item(5, 12).Value = 1.024804589533659
$ws.Cells.Item(5, 13).Value = 1.042290760949904
$ws.Cells.Item(5, 14).Value = 1.012745942495637
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.021547621791594
$ws.Cells.Item(6, 4).Value = 1.030907788868946
$ws.Cells.Item(6, 5).Value = 1.022370620503709
$ws.Cells.Item(6, 6).Value = 1.039908183339641
$ws.Cells.Item(6, 9).Value = 1.029536237740329
$ws.Cells.Item(6, 10).Value = 1.026034486301041
$ws.Cells.Item(6, 11).Value = 1.033343887330546
$ws.Cells.Item(6, 12).Value = 1.024828173430598
$ws.Cells.Item(6, 13).Value = 1.04232207361481
$ws.Cells.Item(6, 14).Value = 1.012753214713622
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.021283584978021
$ws.Cells.Item(7, 4).Value = 1.03070646813696
$ws.Cells.Item(7, 5).Value = 1.022147223933778
$ws.Cells.Item(7, 6).Value = 1.039632661736264
$ws.Cells.Item(7, 9).Value = 1.029505487353829
$ws.Cells.Item(7, 10).Value = 1.025884781360952
$ws.Cells.Item(7, 11).Value = 1.033203198753464
$ws.Cells.Item(7, 12).Value = 1.024666031711279
$ws.Cells.Item(7, 13).Value = 1.042106785978813
$ws.Cells.Item(7, 14).Value = 1.012703207037513
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.020181356633675
$ws.Cells.Item(8, 4).Value = 1.029865785518977
$ws.Cells.Item(8, 5).Value = 1.021215472203944
$ws.Cells.Item(8, 6).Value = 1.038482935369966
$ws.Cells.Item(8, 9).Value = 1.029374434642699
$ws.Cells.Item(8, 10).Value = 1.025259154057704
$ws.Cells.Item(8, 11).Value = 1.03261455536159
$ws.Cells.Item(8, 12).Value = 1.023988953278713
$ws.Cells.Item(8, 13).Value = 1.041207532639477
$ws.Cells.Item(8, 14).Value = 1.012494123272111
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.018244482362803
$ws.Cells.Item(9, 4).Value = 1.028387552069025
$ws.Cells.Item(9, 5).Value = 1.019581301251774
$ws.Cells.Item(9, 6).Value = 1.036464297357059
$ws.Cells.Item(9, 9).Value = 1.029134043960266
$ws.Cells.Item(9, 10).Value = 1.024157235629213
$ws.Cells.Item(9, 11).Value = 1.031575164254302
$ws.Cells.Item(9, 12).Value = 1.0227983859271
$ws.Cells.Item(9, 13).Value = 1.039625361280612
$ws.Cells.Item(9, 14).Value = 1.012125496347768
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.016956954617543
$ws.Cells.Item(10, 4).Value = 1.027404295361449
$ws.Cells.Item(10, 5).Value = 1.01849713575805
$ws.Cells.Item(10, 6).Value = 1.03512358290347
$ws.Cells.Item(10, 9).Value = 1.028967463872077
$ws.Cells.Item(10, 10).Value = 1.023423046345704
$ws.Cells.Item(10, 11).Value = 1.030880889141335
$ws.Cells.Item(10, 12).Value = 1.022006458162339
$ws.Cells.Item(10, 13).Value = 1.038572315738078
$ws.Cells.Item(10, 14).Value = 1.011879642167419
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.01640034379926
$ws.Cells.Item(11, 4).Value = 1.026979087032791
$ws.Cells.Item(11, 5).Value = 1.018028955828821
$ws.Cells.Item(11, 6).Value = 1.034544260089837
$ws.Cells.Item(11, 9).Value = 1.028893842367215
$ws.Cells.Item(11, 10).Value = 1.023105250538981
$ws.Cells.Item(11, 11).Value = 1.030579957965095
$ws.Cells.Item(11, 12).Value = 1.021663985080347
$ws.Cells.Item(11, 13).Value = 1.038116768718095
$ws.Cells.Item(11, 14).Value = 1.011773166084203
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.016193730167454
$ws.Cells.Item(12, 4).Value = 1.026821229985531
$ws.Cells.Item(12, 5).Value = 1.017855245757914
$ws.Cells.Item(12, 6).Value = 1.034329258416296
$ws.Cells.Item(12, 9).Value = 1.028866272685972
$ws.Cells.Item(12, 10).Value = 1.022987225284883
$ws.Cells.Item(12, 11).Value = 1.03046813422126
$ws.Cells.Item(12, 12).Value = 1.021536842360822
$ws.Cells.Item(12, 13).Value = 1.037947624421713
$ws.Cells.Item(12, 14).Value = 1.01173361365592
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.016238043294055
$ws.Cells.Item(13, 4).Value = 1.026855087014997
$ws.Cells.Item(13, 5).Value = 1.017892498403255
$ws.Cells.Item(13, 6).Value = 1.034375368618597
$ws.Cells.Item(13, 9).Value = 1.028872196573773
$ws.Cells.Item(13, 10).Value = 1.023012541255729
$ws.Cells.Item(13, 11).Value = 1.030492122785631
$ws.Cells.Item(13, 12).Value = 1.021564111850003
$ws.Cells.Item(13, 13).Value = 1.037983903425929
$ws.Cells.Item(13, 14).Value = 1.011742097890214
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.016383262259935
$ws.Cells.Item(14, 4).Value = 1.026966036790311
$ws.Cells.Item(14, 5).Value = 1.018014592956325
$ws.Cells.Item(14, 6).Value = 1.034526484206714
$ws.Cells.Item(14, 9).Value = 1.028891568001731
$ws.Cells.Item(14, 10).Value = 1.023095494154782
$ws.Cells.Item(14, 11).Value = 1.030570715479434
$ws.Cells.Item(14, 12).Value = 1.021653474047984
$ws.Cells.Item(14, 13).Value = 1.038102785838229
$ws.Cells.Item(14, 14).Value = 1.01176989671579
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.016472754581593
$ws.Cells.Item(15, 4).Value = 1.027034407829936
$ws.Cells.Item(15, 5).Value = 1.018089845043767
$ws.Cells.Item(15, 6).Value = 1.034619616031384
$ws.Cells.Item(15, 9).Value = 1.028903473798167
$ws.Cells.Item(15, 10).Value = 1.023146606631422
$ws.Cells.Item(15, 11).Value = 1.030619133175978
$ws.Cells.Item(15, 12).Value = 1.021708541971621
$ws.Cells.Item(15, 13).Value = 1.038176042033939
$ws.Cells.Item(15, 14).Value = 1.011787024176576
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.016993914085086
$ws.Cells.Item(16, 4).Value = 1.02743252675076
$ws.Cells.Item(16, 5).Value = 1.01852823425332
$ws.Cells.Item(16, 6).Value = 1.035162056429558
$ws.Cells.Item(16, 9).Value = 1.028972318527079
$ws.Cells.Item(16, 10).Value = 1.023444139896915
$ws.Cells.Item(16, 11).Value = 1.03090085463495
$ws.Cells.Item(16, 12).Value = 1.022029196296471
$ws.Cells.Item(16, 13).Value = 1.038602558074046
$ws.Cells.Item(16, 14).Value = 1.011886708261955
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.017321065016423
$ws.Cells.Item(17, 4).Value = 1.027682404359285
$ws.Cells.Item(17, 5).Value = 1.018803565791142
$ws.Cells.Item(17, 6).Value = 1.035502641777558
$ws.Cells.Item(17, 9).Value = 1.029015104231072
$ws.Cells.Item(17, 10).Value = 1.023630805811332
$ws.Cells.Item(17, 11).Value = 1.031077490259764
$ws.Cells.Item(17, 12).Value = 1.022230452231281
$ws.Cells.Item(17, 13).Value = 1.038870216314844
$ws.Cells.Item(17, 14).Value = 1.011949232542669
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.017511973077142
$ws.Cells.Item(18, 4).Value = 1.027828206519167
$ws.Cells.Item(18, 5).Value = 1.018964284515874
$ws.Cells.Item(18, 6).Value = 1.03570141658138
$ws.Cells.Item(18, 9).Value = 1.029039916481963
$ws.Cells.Item(18, 10).Value = 1.023739695697804
$ws.Cells.Item(18, 11).Value = 1.031180489264771
$ws.Cells.Item(18, 12).Value = 1.022347883470655
$ws.Cells.Item(18, 13).Value = 1.03902637807931
$ws.Cells.Item(18, 14).Value = 1.011985699974793
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.017577082407897
$ws.Cells.Item(19, 4).Value = 1.027877930182708
$ws.Cells.Item(19, 5).Value = 1.019019106168512
$ws.Cells.Item(19, 6).Value = 1.035769213459105
$ws.Cells.Item(19, 9).Value = 1.029048352385719
$ws.Cells.Item(19, 10).Value = 1.023776826130829
$ws.Cells.Item(19, 11).Value = 1.031215604206528
$ws.Cells.Item(19, 12).Value = 1.022387931579748
$ws.Cells.Item(19, 13).Value = 1.039079632147351
$ws.Cells.Item(19, 14).Value = 1.011998134086123
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.017285955863254
$ws.Cells.Item(20, 4).Value = 1.027655589384018
$ws.Cells.Item(20, 5).Value = 1.018774012638958
$ws.Cells.Item(20, 6).Value = 1.035466088062873
$ws.Cells.Item(20, 9).Value = 1.029010528612146
$ws.Cells.Item(20, 10).Value = 1.023610777197228
$ws.Cells.Item(20, 11).Value = 1.031058541977073
$ws.Cells.Item(20, 12).Value = 1.022208855013845
$ws.Cells.Item(20, 13).Value = 1.038841494840314
$ws.Cells.Item(20, 14).Value = 1.011942524474551
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.016340495115412
$ws.Cells.Item(21, 4).Value = 1.0269333625058
$ws.Cells.Item(21, 5).Value = 1.017978633816372
$ws.Cells.Item(21, 6).Value = 1.034481979320044
$ws.Cells.Item(21, 9).Value = 1.028885869759965
$ws.Cells.Item(21, 10).Value = 1.023071066079246
$ws.Cells.Item(21, 11).Value = 1.030547573107426
$ws.Cells.Item(21, 12).Value = 1.02162715724635
$ws.Cells.Item(21, 13).Value = 1.03806777609075
$ws.Cells.Item(21, 14).Value = 1.011761710718904
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.015746836442169
$ws.Cells.Item(22, 4).Value = 1.026479758224139
$ws.Cells.Item(22, 5).Value = 1.017479663987628
$ws.Cells.Item(22, 6).Value = 1.033864300210203
$ws.Cells.Item(22, 9).Value = 1.028806199734337
$ws.Cells.Item(22, 10).Value = 1.022731834549804
$ws.Cells.Item(22, 11).Value = 1.030226049886843
$ws.Cells.Item(22, 12).Value = 1.02126180928832
$ws.Cells.Item(22, 13).Value = 1.037581691994325
$ws.Cells.Item(22, 14).Value = 1.01164801187742
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.016061470711486
$ws.Cells.Item(23, 4).Value = 1.026720175512482
$ws.Cells.Item(23, 5).Value = 1.01774407089037
$ws.Cells.Item(23, 6).Value = 1.034191641580542
$ws.Cells.Item(23, 9).Value = 1.02884855657265
$ws.Cells.Item(23, 10).Value = 1.022911657100851
$ws.Cells.Item(23, 11).Value = 1.030396519309237
$ws.Cells.Item(23, 12).Value = 1.021455449860317
$ws.Cells.Item(23, 13).Value = 1.037839337426502
$ws.Cells.Item(23, 14).Value = 1.011708286958614
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.017301819911368
$ws.Cells.Item(24, 4).Value = 1.027667705754268
$ws.Cells.Item(24, 5).Value = 1.018787366058077
$ws.Cells.Item(24, 6).Value = 1.035482604751431
$ws.Cells.Item(24, 9).Value = 1.029012596582171
$ws.Cells.Item(24, 10).Value = 1.023619827232557
$ws.Cells.Item(24, 11).Value = 1.031067103981928
$ws.Cells.Item(24, 12).Value = 1.022218613736989
$ws.Cells.Item(24, 13).Value = 1.038854472711138
$ws.Cells.Item(24, 14).Value = 1.011945555567789
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.018744561833891
$ws.Cells.Item(25, 4).Value = 1.028769325430647
$ws.Cells.Item(25, 5).Value = 1.020002850585036
$ws.Cells.Item(25, 6).Value = 1.036985283241588
$ws.Cells.Item(25, 9).Value = 1.029197307531388
$ws.Cells.Item(25, 10).Value = 1.024442038927174
$ws.Cells.Item(25, 11).Value = 1.031844115302529
$ws.Cells.Item(25, 12).Value = 1.023105867707583
$ws.Cells.Item(25, 13).Value = 1.040034092170102
$ws.Cells.Item(25, 14).Value = 1.012220815424692
